# Wales Premier League workbook update
# Swaps the content of several fixture rows (columns B:AD) so that each
# listed row receives the data that another (paired) row held before the
# edit. Column A (the running row index) is left untouched.
#
# This reproduces the diff where, for each pair/triple of rows below,
# the full fixture record (id, date, teams, scores, odds, ...) moves to
# a different row position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Swap-Rows {
    param(
        [int]$RowA,
        [int]$RowB
    )
    $addrA = "B" + $RowA + ":AD" + $RowA
    $addrB = "B" + $RowB + ":AD" + $RowB
    $rangeA = $ws.Range($addrA)
    $rangeB = $ws.Range($addrB)

    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2

    $rangeA.Value2 = $valuesB
    $rangeB.Value2 = $valuesA
}

function Rotate-Rows {
    # Row $Rows[0] receives the data that was in $Rows[1],
    # row $Rows[1] receives the data that was in $Rows[2],
    # ...
    # row $Rows[-1] receives the data that was in $Rows[0].
    param(
        [int[]]$Rows
    )
    $n = $Rows.Length
    $ranges = @()
    $values = @()
    for ($i = 0; $i -lt $n; $i++) {
        $r = $Rows[$i]
        $addr = "B" + $r + ":AD" + $r
        $rng = $ws.Range($addr)
        $ranges += $rng
        $values += ,($rng.Value2)
    }
    for ($i = 0; $i -lt $n; $i++) {
        $srcIndex = ($i + 1) % $n
        $ranges[$i].Value2 = $values[$srcIndex]
    }
}

# Simple pairwise row swaps
Swap-Rows 8   10
Swap-Rows 14  15
Swap-Rows 20  23
Swap-Rows 21  22
Swap-Rows 26  28
Swap-Rows 56  57
Swap-Rows 59  60
Swap-Rows 68  69
Swap-Rows 80  81
Swap-Rows 163 164
Swap-Rows 170 172
Swap-Rows 175 176

# Three-way rotation: 128 <- 129 <- 130 <- 128
Rotate-Rows @(128, 129, 130)

$wb.Save()
